$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "general": summary stats for this run
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value  = 145.1089212663341      # objValue
$ws.Range("B4").Value  = 0.01399993896484375    # runtime
$ws.Range("B6").Value  = 37.36892126633407      # Z1
$ws.Range("B9").Value  = 0                      # Z4
$ws.Range("B10").Value = 107.74                 # Z5

# ---------------------------------------------------------------
# Sheet "x": (i, j, x) assignment table
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("x")
$ws.Range("B4").Value  = 3
$ws.Range("B6").Value  = 9
$ws.Range("B7").Value  = 6
$ws.Range("B10").Value = 4
$ws.Range("B12").Value = 13
$ws.Range("B14").Value = 11

# ---------------------------------------------------------------
# Sheet "U"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("U")
$ws.Range("B4").Value = 2
$ws.Range("B7").Value = 2

# ---------------------------------------------------------------
# Sheet "TBar"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B5").Value  = 10
$ws.Range("B7").Value  = 26.09699127290008
$ws.Range("B8").Value  = 10.34885527085025
$ws.Range("B11").Value = 24.14711948224307
$ws.Range("B13").Value = 26.71671453559703
$ws.Range("B14").Value = 27.28210730097497
$ws.Range("B15").Value = 27.52855135814675

# ---------------------------------------------------------------
# Sheet "Q"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Q")
$ws.Range("C12").Value = 81.47500000000072
$ws.Range("C13").Value = 80.68000000000072
$ws.Range("C16").Value = 87.34500000000074
$ws.Range("C17").Value = 40.35
$ws.Range("C18").Value = 30.90499999999942
$ws.Range("C19").Value = 27.59499999999942
$ws.Range("C20").Value = 31.97499999999942
$ws.Range("C21").Value = 33.99499999999941
$ws.Range("C22").Value = 216.7599999999987
$ws.Range("C23").Value = 216.7249999999988
$ws.Range("C24").Value = 199.5399999999988
$ws.Range("C25").Value = 218.1849999999988
$ws.Range("C26").Value = 206.5849999999988
$ws.Range("C27").Value = 236.8400000000015
$ws.Range("C28").Value = 244.9650000000014
$ws.Range("C29").Value = 236.8400000000015
$ws.Range("C30").Value = 243.8100000000014
$ws.Range("C31").Value = 244.9650000000015
$ws.Range("C32").Value = 144.1299999999993
$ws.Range("C33").Value = 136.6699999999993
$ws.Range("C34").Value = 117.2199999999993
$ws.Range("C35").Value = 134.2299999999993
$ws.Range("C36").Value = 125.7699999999993
$ws.Range("C42").Value = 284.6849999999982
$ws.Range("C43").Value = 295.8849999999982
$ws.Range("C44").Value = 259.3599999999982
$ws.Range("C45").Value = 281.9549999999982
$ws.Range("C46").Value = 265.5149999999982
$ws.Range("C47").Value = 166.9600000000012
$ws.Range("C48").Value = 168.6450000000012
$ws.Range("C49").Value = 164.4300000000013
$ws.Range("C50").Value = 171.1650000000012
$ws.Range("C51").Value = 172.0750000000012
$ws.Range("C53").Value = 260.9900000000009
$ws.Range("C54").Value = 252.975000000001
$ws.Range("C57").Value = 250.970000000001
$ws.Range("C58").Value = 260.9900000000009
$ws.Range("C59").Value = 252.975000000001
$ws.Range("C60").Value = 269.580000000001
$ws.Range("C61").Value = 250.575000000001
$ws.Range("C62").Value = 236.8400000000015
$ws.Range("C63").Value = 244.9650000000014
$ws.Range("C64").Value = 236.8400000000015
$ws.Range("C65").Value = 243.8100000000014
$ws.Range("C66").Value = 244.9650000000015
$ws.Range("C67").Value = 284.6849999999982
$ws.Range("C68").Value = 295.8849999999982
$ws.Range("C69").Value = 259.3599999999982
$ws.Range("C70").Value = 281.9549999999982
$ws.Range("C71").Value = 265.5149999999982

# ---------------------------------------------------------------
# Sheet "R"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("R")
$ws.Range("C8").Value  = 0
$ws.Range("C10").Value = 0

# ---------------------------------------------------------------
# Sheet "L"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("L")
$ws.Range("C12").Value = 6.57
$ws.Range("C13").Value = 5.2
$ws.Range("C14").Value = 7.32
$ws.Range("C15").Value = 5.515
$ws.Range("C16").Value = 5.44
$ws.Range("C27").Value = 3.6
$ws.Range("C28").Value = 6.475
$ws.Range("C29").Value = 4.165
$ws.Range("C30").Value = 6.58
$ws.Range("C31").Value = 3.005

# ---------------------------------------------------------------
# Sheet "rho": drop the last two data rows (A2:C3)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("rho")
$ws.Range("A2:C3").EntireRow.Delete()
